# Scheduled-runner update: refresh computed market-profit columns (H:N)
# on the Leve profitability sheets, per latest price pull.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 241.11111
$ws.Range("I2").Value = 195
$ws.Range("K2").Value = 195
$ws.Range("M2").Value = -82
$ws.Range("H61").Value = 1512.4286
$ws.Range("I61").Value = 1512.4286
$ws.Range("K61").Value = 4537.2858
$ws.Range("M61").Value = -4365.2858
$ws.Range("H62").Value = 9195.429
$ws.Range("I62").Value = 8097
$ws.Range("J62").Value = 10660
$ws.Range("K62").Value = 8097
$ws.Range("L62").Value = 10660
$ws.Range("M62").Value = -7473
$ws.Range("N62").Value = -11908
$ws.Range("H65").Value = 9195.429
$ws.Range("I65").Value = 8097
$ws.Range("J65").Value = 10660
$ws.Range("K65").Value = 40485
$ws.Range("L65").Value = 53300
$ws.Range("M65").Value = -37365
$ws.Range("N65").Value = -59540
$ws.Range("H135").Value = 566.2222
$ws.Range("I135").Value = 497
$ws.Range("K135").Value = 4473
$ws.Range("M135").Value = -1938
$ws.Range("H137").Value = 2425
$ws.Range("I137").Value = 2186.8235
$ws.Range("J137").Value = 2874.889
$ws.Range("K137").Value = 6560.470499999999
$ws.Range("L137").Value = 8624.667000000001
$ws.Range("M137").Value = -4010.470499999999
$ws.Range("N137").Value = -13724.667
$ws.Range("H141").Value = 1039313.94
$ws.Range("I141").Value = 1274662.6
$ws.Range("J141").Value = 3780
$ws.Range("K141").Value = 3823987.8
$ws.Range("L141").Value = 11340
$ws.Range("M141").Value = -3818807.8
$ws.Range("N141").Value = -21700

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3464.6943
$ws.Range("I32").Value = 3000.7856
$ws.Range("K32").Value = 3000.7856
$ws.Range("M32").Value = -2713.7856
$ws.Range("H61").Value = 3236.5
$ws.Range("I61").Value = 1529.8889
$ws.Range("K61").Value = 1529.8889
$ws.Range("M61").Value = -1317.8889
$ws.Range("H74").Value = 999.5599999999999
$ws.Range("J74").Value = 2775.25
$ws.Range("L74").Value = 2775.25
$ws.Range("N74").Value = -4523.25
$ws.Range("H77").Value = 999.5599999999999
$ws.Range("J77").Value = 2775.25
$ws.Range("L77").Value = 13876.25
$ws.Range("N77").Value = -22612.25
$ws.Range("H136").Value = 3236.5
$ws.Range("I136").Value = 1529.8889
$ws.Range("K136").Value = 4589.6667
$ws.Range("M136").Value = -2039.6667

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1027.8334
$ws.Range("I99").Value = 1027.8334
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1027.8334
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 470.1666
$ws.Range("N99").ClearContents()
$ws.Range("H134").Value = 4159.0435
$ws.Range("I134").Value = 4159.0435
$ws.Range("K134").Value = 12477.1305
$ws.Range("M134").Value = -9942.130499999999

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1360.3226
$ws.Range("J31").Value = 2140.3076
$ws.Range("L31").Value = 2140.3076
$ws.Range("N31").Value = -2730.3076
$ws.Range("H34").Value = 1360.3226
$ws.Range("J34").Value = 2140.3076
$ws.Range("L34").Value = 2140.3076
$ws.Range("N34").Value = -2544.3076
$ws.Range("H58").Value = 3345354.2
$ws.Range("I58").Value = 3953483.2
$ws.Range("J58").Value = 644
$ws.Range("K58").Value = 3953483.2
$ws.Range("L58").Value = 644
$ws.Range("M58").Value = -3953280.2
$ws.Range("N58").Value = -1050
$ws.Range("H107").Value = 533.4167
$ws.Range("I107").Value = 440.46667
$ws.Range("J107").Value = 688.3333
$ws.Range("K107").Value = 440.46667
$ws.Range("L107").Value = 688.3333
$ws.Range("M107").Value = 1479.53333
$ws.Range("N107").Value = -4528.3333
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()
$ws.Range("H132").Value = 2470.682
$ws.Range("I132").Value = 1644.75
$ws.Range("J132").Value = 4673.1665
$ws.Range("K132").Value = 4934.25
$ws.Range("L132").Value = 14019.4995
$ws.Range("M132").Value = -2404.25
$ws.Range("N132").Value = -19079.4995
$ws.Range("H136").Value = 3345354.2
$ws.Range("I136").Value = 3953483.2
$ws.Range("J136").Value = 644
$ws.Range("K136").Value = 11860449.6
$ws.Range("L136").Value = 1932
$ws.Range("M136").Value = -11857899.6
$ws.Range("N136").Value = -7032

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 33694.637
$ws.Range("J129").Value = 56566.383
$ws.Range("L129").Value = 169699.149
$ws.Range("N129").Value = -179699.149
$ws.Range("H131").Value = 11645835
$ws.Range("J131").Value = 19223.625
$ws.Range("L131").Value = 57670.875
$ws.Range("N131").Value = -67750.875

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3005
$ws.Range("I80").Value = 3005
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 3005
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -2007
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 3005
$ws.Range("I83").Value = 3005
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 15025
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -10033
$ws.Range("N83").Value = 0
$ws.Range("H107").Value = 1050.75
$ws.Range("I107").Value = 100
$ws.Range("J107").Value = 1367.6666
$ws.Range("K107").Value = 100
$ws.Range("L107").Value = 1367.6666
$ws.Range("M107").Value = 1820
$ws.Range("N107").Value = -5207.6666
$ws.Range("H126").Value = 2528053.5
$ws.Range("I126").Value = 3270845
$ws.Range("J126").Value = 2562.6
$ws.Range("K126").Value = 9812535
$ws.Range("L126").Value = 7687.799999999999
$ws.Range("M126").Value = -9810065
$ws.Range("N126").Value = -12627.8

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2000
$ws.Range("J22").Value = 2000
$ws.Range("L22").Value = 2000
$ws.Range("N22").Value = -2590
$ws.Range("H27").Value = 2000
$ws.Range("J27").Value = 2000
$ws.Range("L27").Value = 2000
$ws.Range("N27").Value = -2214
$ws.Range("H40").Value = 5727.2856
$ws.Range("I40").Value = 4810.625
$ws.Range("K40").Value = 4810.625
$ws.Range("M40").Value = -4674.625
$ws.Range("H132").Value = 6363.8125
$ws.Range("I132").Value = 1101
$ws.Range("K132").Value = 3303
$ws.Range("M132").Value = -773
$ws.Range("H136").Value = 3876.5
$ws.Range("I136").Value = 2166.2
$ws.Range("K136").Value = 6498.599999999999
$ws.Range("M136").Value = -3948.599999999999
